$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-39 down to 27-40
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new data record
$ws.Cells.Item(26, 1).Value = 7
$ws.Cells.Item(26, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value = "Ñuble"
$ws.Cells.Item(26, 4).Value = 44893
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat
$ws.Cells.Item(26, 5).Value = 16
$ws.Cells.Item(26, 6).Value = 300000000
$ws.Cells.Item(26, 7).Value = "Espárragos"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 2000
$ws.Cells.Item(26, 11).Value = 900
$ws.Cells.Item(26, 12).Value = 1000
$ws.Cells.Item(26, 13).Value = 950
$ws.Cells.Item(26, 14).Value = "$/kilo"
$ws.Cells.Item(26, 15).Value = "Región de Ñuble"
$ws.Cells.Item(26, 16).Value = 950
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = "Hortaliza"
